$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.661.48"
$ws.Range("E2").Value = "  +4.75%  "

$ws.Range("D3").Value = "3.333.44"
$ws.Range("E3").Value = "  +4.51%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").Value = "3.905.11"
$ws.Range("E12").Value = "  +4.44%  "

$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("D16").Value = "62.705.96"
$ws.Range("E16").Value = "  +4.70%  "

$ws.Range("D17").Value = "3.331.05"
$ws.Range("E17").Value = "  +4.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.09%  "

$ws.Range("E19").Value = "  +5.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.538"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.177"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "

$ws.Range("D27").Value = "0.0₃0966"
$ws.Range("E27").Value = "  +6.36%  "

$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.18%  "

$ws.Range("E35").Value = "  +9.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("E37").Value = "  +11.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.28%  "

$ws.Range("D39").Value = "2.863.05"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0738"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.95%  "

$ws.Range("E41").Value = "  +8.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.753"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.98%  "

$ws.Range("B47").Value = "RenzoRestakedETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D47").Value = "3.376.60"
$ws.Range("E47").Value = "  +4.30%  "

$ws.Range("E48").Value = "  +3.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.805"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "284.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.08%  "
